# Loan RBI, Variable Instalments
# Switch focus to the "Repayment schedule" sheet and insert a new
# (currently blank) column before the old "Late" column so that a new
# data column can be tracked going forward.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make this the active/visible sheet (was "NewLoanInput" before).
$ws.Activate()

# Insert a new blank column before column N ("Late"); this pushes the
# existing Late / heading / Outstanding columns one to the right.
$ws.Columns("N:N").Insert()

# Give the freshly inserted column a fixed width (matches the width
# used for the neighbouring "In Advance" column).
$ws.Columns("N:N").ColumnWidth = 10.14

# Leave the selection on L12, as last used on this sheet.
$ws.Range("L12").Select()
